$wb = $excel.ActiveWorkbook

# --- Sheet1: update the two step-label groups in column A ---
$ws1 = $wb.Worksheets.Item("Sheet1")

# Rows 5-14 previously read "Upload Data File - {Org, Submission}";
# relabel them as "Step 1 ..." per the commit.
$ws1.Range("A5:A14").Value = "Step 1 Upload Data File - {Org, Submission}"

# Rows 16-24 previously read "Validate Data File - {Org, Submission}";
# relabel them as "Step 2 ..." per the commit.
$ws1.Range("A16:A24").Value = "Step 2 Validate Data File - {Org, Submission}"

# Update the saved selection/view on Sheet1.
$ws1.Activate()
$ws1.Range("F23").Select()

# --- CycleGroup1 sheet: update the saved selection ---
$ws2 = $wb.Worksheets.Item("CycleGroup1")
$ws2.Activate()
$ws2.Range("C33").Select()

# Leave the final active sheet as Sheet1 (it was the tab-selected one).
$ws1.Activate()
